# Build site at 2023-04-12 14:53:07 UTC
# Updates the LOT2052 course sheet: fills in previously-empty
# "Objetivos / Programa resumido / Programa / Bibliografia" content,
# inserts a new row so the "Docentes responsaveis" value lines up with its
# own row, and shifts the "Requisitos" list down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row at position 13 ("Programa resumido:" and everything
#    below shifts down by one row), matching the new dimension A1:C25.
# ---------------------------------------------------------------------
$ws.Range("A13").EntireRow.Insert()

# The inserted row copies the formatting of the row above it, leaving a
# stray formatted-but-empty A13 cell. The target layout has no cell in
# column A on this row at all, so drop it.
$ws.Range("A13").Clear()

# ---------------------------------------------------------------------
# 2) Fill the new row 13 with the "Docentes responsaveis" value (moved
#    down from the old row 10) in columns B and C.
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "1097178 - João Batista de Almeida e Silva"
$ws.Range("C13").Value = "1097178 - João Batista de Almeida e Silva"

# B13 needs to match the normal "value" cell style (non-bold, wrap text,
# top-aligned) used throughout column B, same as C13's column style.
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160

# ---------------------------------------------------------------------
# 3) Fill in the previously-missing content cells.
# ---------------------------------------------------------------------

# Objetivos / Objectives text
$ws.Range("B10").Value = "Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas."
$ws.Range("C10").Value = "Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas."

# Programa resumido (short syllabus, pt-br)
$ws.Range("B14").Value = "Elaboração prática de cerveja, cachaça, fermentados e destilados de frutas, cereais e tuberculos, vinhos e análise sensorial."
$ws.Range("C14").Value = "Elaboração prática de cerveja, cachaça, fermentados e destilados de frutas, cereais e tuberculos, vinhos e análise sensorial."

# Programa (full syllabus, pt-br)
$ws.Range("B16").Value = "1. Elaboração de cerveja: matérias-primas, preparação do mosto, tecnologia de fermentação e maturação.2. Elaboração de aguardente: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.3. Elaboração de destilados de frutas: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.4. Elaboração e vinhos: matérias-primas, preparação do mosto, tecnologia de fermentação, maturação.5. Análise sensorial: teste sensorial das bebidas preparadas nos itens anteriores"
$ws.Range("C16").Value = "1. Elaboração de cerveja: matérias-primas, preparação do mosto, tecnologia de fermentação e maturação.2. Elaboração de aguardente: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.3. Elaboração de destilados de frutas: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.4. Elaboração e vinhos: matérias-primas, preparação do mosto, tecnologia de fermentação, maturação.5. Análise sensorial: teste sensorial das bebidas preparadas nos itens anteriores"

# Método (evaluation method)
$ws.Range("B19").Value = "Relatórios e seminários sobre os experimentos"
$ws.Range("C19").Value = "Relatórios e seminários sobre os experimentos"

# Critério (evaluation criteria)
$ws.Range("B20").Value = "Média aritmética entre os relatórios e seminários"
$ws.Range("C20").Value = "Média aritmética entre os relatórios e seminários"

# Norma de recuperação (recovery policy)
$ws.Range("B21").Value = "A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2."
$ws.Range("C21").Value = "A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2."

# Bibliografia
$ws.Range("B22").Value = "1. AQUARONE, E.; BORZANI, W.; SCHMIDELL, W.; LIMA, U. A. Biotecnologia na Produção deAlimentos. V. 4, Biotecnologia Industrial, São Paulo: Edgard Blücher Ltda. 2001.4. DUVAL, G. Fabricação de Vinhos de Frutas. S.I.A.RJ:Ministério da Agricultura, 1947.5. HOUGH, J.S. Biotecnología de La cerveza y de la malta. Editorial ACRIBA S/A, 1978.6. LIMA, U. A. Aguardente: fabricação em pequenas destilarias. Ed. FEALQ. 1999.7. MARTINELLI FILHO, A. Tecnologia de Vinhos e Vinagres de Frutas. Agroindústria de BaixoInvestimento. Departamento de Tecnologia Rural da ESALQ/USP.8. MORRETO, E. et al. Vinhos e Vinagres: Processamento e Análises. FlorianópolisEditoraUFSC, 1988.9. PACHECO, A. O. Manual do Bar. São Paulo. Editora SENAC, 1996.10. STANIER, R. Y.; INGRAHAM, J. L., WHEELIS, M. L.; PAINTER, P. R. The Microbial World.Englewood Cliffs, New Jersey, 1986.11.Venturini Filho, W.G. Bebidas Alcoólicas. Ciência e Tecnologia. São Paulo. Edgar Blucher Ltda. 2a. Edição. 2016. 575 p."
$ws.Range("C22").Value = "1. AQUARONE, E.; BORZANI, W.; SCHMIDELL, W.; LIMA, U. A. Biotecnologia na Produção deAlimentos. V. 4, Biotecnologia Industrial, São Paulo: Edgard Blücher Ltda. 2001.4. DUVAL, G. Fabricação de Vinhos de Frutas. S.I.A.RJ:Ministério da Agricultura, 1947.5. HOUGH, J.S. Biotecnología de La cerveza y de la malta. Editorial ACRIBA S/A, 1978.6. LIMA, U. A. Aguardente: fabricação em pequenas destilarias. Ed. FEALQ. 1999.7. MARTINELLI FILHO, A. Tecnologia de Vinhos e Vinagres de Frutas. Agroindústria de BaixoInvestimento. Departamento de Tecnologia Rural da ESALQ/USP.8. MORRETO, E. et al. Vinhos e Vinagres: Processamento e Análises. FlorianópolisEditoraUFSC, 1988.9. PACHECO, A. O. Manual do Bar. São Paulo. Editora SENAC, 1996.10. STANIER, R. Y.; INGRAHAM, J. L., WHEELIS, M. L.; PAINTER, P. R. The Microbial World.Englewood Cliffs, New Jersey, 1986.11.Venturini Filho, W.G. Bebidas Alcoólicas. Ciência e Tecnologia. São Paulo. Edgar Blucher Ltda. 2a. Edição. 2016. 575 p."

# ---------------------------------------------------------------------
# 4) Column layout cleanup: column A's width entry used to span columns
#    A:B (max="2"); narrow it back to column A only so it no longer
#    overlaps column B's own (wider) width entry.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).Hidden = $false
